# Updated cryptos list on Sat Feb  3 23:56:22 UTC 2024 with GitHub Actions
# Refreshes Price/Volume(1h) figures and inserts a new "LEO" row at position 27,
# which pushes the existing rows 27-50 down by one (dropping "Stacks" off the
# bottom of the 50-row table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.980.38"
$ws.Range("D3").Value = "2.295.10"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "299.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.63%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.519"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.513"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0786"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("E12").Value = "  +0.80%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.66"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.41%  "
$ws.Range("D15").Value = "2.651.04"
$ws.Range("E15").Value = "  -0.55%  "
$ws.Range("D16").Value = "2.309.51"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.785"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "42.896.56"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("E21").Value = "  -0.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.90%  "
$ws.Range("E24").Value = "  -2.78%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("E26").Value = "  -0.96%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.19%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.06%  "
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "165.10"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.99%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.04"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.05"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.18%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.61%  "
$ws.Range("B33").Value = "FirstDigitalUSD"
$ws.Range("C33").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.06"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.71"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.67%  "
$ws.Range("B36").Value = "Celestia"
$ws.Range("C36").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.74"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.41"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0693"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.13%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.102"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.77"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.44%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.111"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("B42").Value = "LidoDAOToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.86%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "2.002.95"
$ws.Range("E43").Value = "  +1.13%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0285"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.63%  "
$ws.Range("B45").Value = "ApeXProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.23"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.65%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.26%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.94%  "
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.42%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.518.27"
$ws.Range("E50").Value = "  -0.51%  "
$ws.Range("B51").Value = "BitcoinSV"
$ws.Range("C51").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.64%  "
